# common ssaDMPs and ssaVMPs
# Replace the existing cg/gene rows with the new shared set and add a new
# "reverses" (yes/no) column E.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data for columns A (probe id) and B (gene, optional) -------------
$rows = @(
    @{ A = "cg15634877"; B = "SPEG";  E = "yes" },
    @{ A = "cg10317815"; B = $null;   E = "yes" },
    @{ A = "cg08464177"; B = "BAI1";  E = "yes" },
    @{ A = "cg16786971"; B = $null;   E = "no"  },
    @{ A = "cg15122985"; B = "TCEA2"; E = "no"  },
    @{ A = "cg14506657"; B = $null;   E = "no"  },
    @{ A = "cg24482850"; B = "NUAK2"; E = "no"  },
    @{ A = "cg20403557"; B = "LVRN";  E = "no"  },
    @{ A = "cg06904667"; B = $null;   E = "no"  },
    @{ A = "cg04202002"; B = $null;   E = "no"  }
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $data = $rows[$i]

    $ws.Cells.Item($r, 1).Value = $data.A

    if ($data.B) {
        $ws.Cells.Item($r, 2).Value = $data.B
    } else {
        $ws.Cells.Item($r, 2).Value = ""
    }

    $ws.Cells.Item($r, 3).Value = "auto"
    $ws.Cells.Item($r, 4).Value = "auto"
    $ws.Cells.Item($r, 5).Value = $data.E
}

# --- New header cell E1 ("reverses"), styled like the other headers -------
$ws.Range("E1").Value = "reverses"

# Start from the same look as the existing bold / centered / top-aligned
# header cells (A1:D1), then narrow the border down to left+right only.
$ws.Range("A1").Copy()
$ws.Range("E1").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("E1").Borders.Item(8).LineStyle = -4142  # xlEdgeTop   -> none
$ws.Range("E1").Borders.Item(9).LineStyle = -4142  # xlEdgeBottom -> none
$ws.Range("E1").Borders.Item(7).LineStyle = 1      # xlEdgeLeft  -> thin
$ws.Range("E1").Borders.Item(7).Weight = 2
$ws.Range("E1").Borders.Item(7).Color = 0
$ws.Range("E1").Borders.Item(10).LineStyle = 1     # xlEdgeRight -> thin
$ws.Range("E1").Borders.Item(10).Weight = 2
$ws.Range("E1").Borders.Item(10).Color = 0

$ws.Range("E1").HorizontalAlignment = -4108  # xlCenter
$ws.Range("E1").VerticalAlignment = -4160    # xlTop
$ws.Range("E1").Font.Bold = $true

# --- Misc view bits that moved along with the edit -------------------------
$ws.Application.ActiveWindow.RangeSelection
$ws.Range("G29").Select()
